$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Remove NOT NULL (REQUIRED = "Yes") constraint from the pure FK rows,
# setting them to "No" instead.
$ws.Range("G3").Value = "No"
$ws.Range("G7").Value = "No"
$ws.Range("G8").Value = "No"
$ws.Range("G13").Value = "No"
$ws.Range("G25").Value = "No"

# Update the view/selection state.
$ws.Range("G38").Select()
